# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (cloned layout from "2021-Q4") between
#   "2021-Q4" and "总计", populated with the new quarter's fund holdings.
# - Inserts a new top data row in "总计" summarising the 2022-Q1 quarter
#   (7 holdings, 3.61 亿元), pushing the older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4"
#    (i.e. right before "总计", matching the diff's sheet order).
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $srcSheet)
$newSheet.Name = "2022-Q1"

# Match the page-margin conventions used by the other sheets in this
# workbook (0.75in/1in/0.5in instead of Excel's 0.7in/0.75in/0.3in default).
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# Clone the header/index-column formatting (bold, centered, bordered style)
# from the "2021-Q4" sheet so the new sheet matches the look of its siblings.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$srcSheet.Range("A2:A8").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Numeric index column (A2:A8) — plain integers, same as sibling sheets.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4
$newSheet.Range("A7").Value = 5
$newSheet.Range("A8").Value = 6

$data = @(
    @("001481", "华宝油气(QDII)美元", "39.80", "94.60", "2.24", "0.8915", 8),
    @("162411", "华宝油气(QDII)人民币A", "39.80", "94.60", "2.24", "0.8915", 8),
    @("006679", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A", "14.75", "83.19", "3.94", "0.5812", 5),
    @("162719", "广发道琼斯美国石油开发与生产指数（QDII-LOF）A", "14.75", "83.19", "3.94", "0.5812", 5),
    @("007844", "华宝油气(QDII)人民币C", "12.98", "94.60", "2.24", "0.2908", 8),
    @("006680", "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C", "4.73", "83.19", "3.94", "0.1864", 5),
    @("004243", "广发道琼斯美国石油开发与生产指数（QDII-LOF）C", "4.73", "83.19", "3.94", "0.1864", 5)
)

# Columns B and D:G hold numeric-looking text (fund codes with leading
# zeros, and decimal figures that must keep their trailing zeros, e.g.
# "39.80") — force text formatting before writing so Excel doesn't coerce
# them to numbers.
$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$row = 2
foreach ($item in $data) {
    $newSheet.Cells.Item($row, 2).Value = $item[0]
    $newSheet.Cells.Item($row, 3).Value = $item[1]
    $newSheet.Cells.Item($row, 4).Value = $item[2]
    $newSheet.Cells.Item($row, 5).Value = $item[3]
    $newSheet.Cells.Item($row, 6).Value = $item[4]
    $newSheet.Cells.Item($row, 7).Value = $item[5]
    $newSheet.Cells.Item($row, 8).Value = $item[6]
    $row = $row + 1
}

# The "@" text number-format above leaves a permanent custom style on the
# cells; the source workbook instead leaves these data cells with no
# explicit style at all. Strip it back off (while keeping the stored
# values as text) by painting over with the formatting of a pristine,
# never-touched cell.
$newSheet.Range("Z100").Copy()
$newSheet.Range("B2:B8").PasteSpecial(-4122)   # xlPasteFormats
$newSheet.Range("D2:G8").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new top row for 2022-Q1 and
#    push the existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Give the new A5 index cell the same formatting as the existing index
# column (A2:A4) before writing the shifted-down data into it.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 3.61

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 7
$totalSheet.Range("D3").Value = 2.17

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.45

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 0.53
